# Fix header labels and add a new "PO Forecast" sheet with forecast data.

$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on the existing sheets ---------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" worksheet after "Monthly Trend" ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row (reuse the existing bold/bordered header style) -----
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1").PasteSpecial(-4122)
$wsForecast.Range("A1").Value = "ds"

$wsWeekly.Range("A1").Copy()
$wsForecast.Range("B1").PasteSpecial(-4122)
$wsForecast.Range("B1").Value = "PO_Forecast"

$wsWeekly.Range("A1").Copy()
$wsForecast.Range("C1").PasteSpecial(-4122)
$wsForecast.Range("C1").Value = "yhat_lower"

$wsWeekly.Range("A1").Copy()
$wsForecast.Range("D1").PasteSpecial(-4122)
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Forecast data rows ----------------------------------------------
$data = @(
    @(45298.99999999999, 45, -10.04141555922648, 94.89463786346471),
    @(45305.99999999999, 44, -5.504501601147632, 94.75094240637128),
    @(45319.99999999999, 43, -9.718477246528742, 93.17068477301973),
    @(45326.99999999999, 43, -13.15276842083295, 95.37212967009225),
    @(45333.99999999999, 42, -7.332761515784829, 95.08556305402827),
    @(45340.99999999999, 42, -8.409449124518142, 95.7710665832185),
    @(45361.99999999999, 40, -8.550586588521046, 94.15569614795339),
    @(45389.99999999999, 38, -6.192809823371654, 90.62564415801752),
    @(45396.99999999999, 38, -11.42161825488506, 91.21268816870551),
    @(45403.99999999999, 37, -15.4314124584812, 84.2143259788097),
    @(45410.99999999999, 36, -12.70541625179816, 91.54292080850638),
    @(45438.99999999999, 34, -19.44637413485728, 84.84909415473095),
    @(45445.99999999999, 34, -17.89100328643057, 80.59848363195731),
    @(45452.99999999999, 33, -18.58192145751834, 87.34580422128074),
    @(45459.99999999999, 33, -15.88571611560051, 86.17075773754024),
    @(45473.99999999999, 32, -19.15308617166201, 85.97924587308194),
    @(45487.99999999999, 31, -20.76880624442041, 83.54310726343192),
    @(45592.99999999999, 23, -25.87360281989473, 77.49945014198563),
    @(45599.99999999999, 22, -38.45585193628231, 76.50335660548777),
    @(45606.99999999999, 22, -32.50604405668875, 70.63683375422127),
    @(45613.99999999999, 21, -32.66613310934194, 74.9137140187202),
    @(45620.99999999999, 21, -27.56015455405065, 72.87064095799525),
    @(45627.99999999999, 20, -29.89117076722463, 70.45270571928397),
    @(45634.99999999999, 20, -26.76880738189304, 73.9507813634668),
    @(45641.99999999999, 19, -31.2394852069691, 71.39431982806035),
    @(45648.99999999999, 19, -34.16106523032527, 70.98451437229643),
    @(45655.99999999999, 18, -32.21944509803892, 70.66988438133897)
)

$rowIndex = 2
foreach ($row in $data) {
    # Column A: date-serial, reuse the existing date number-format style.
    $wsWeekly.Range("A2").Copy()
    $wsForecast.Cells.Item($rowIndex, 1).PasteSpecial(-4122)
    $wsForecast.Cells.Item($rowIndex, 1).Value = $row[0]

    $wsForecast.Cells.Item($rowIndex, 2).Value = $row[1]
    $wsForecast.Cells.Item($rowIndex, 3).Value = $row[2]
    $wsForecast.Cells.Item($rowIndex, 4).Value = $row[3]

    $rowIndex++
}

Write-Output "PO Forecast sheet created with $($data.Count) data rows."
